$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.996.76"
$ws.Range("D3").Value = "1.859.37"
$ws.Range("E3").Value = "  -0.86%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'312.43"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").Value = "'0.5140"
$ws.Range("E7").Value = "  +1.31%  "
$ws.Range("D8").Value = "'0.3833"
$ws.Range("E8").Value = "  -0.24%  "
$ws.Range("D9").Value = "'0.08252"
$ws.Range("E9").Value = "  -8.06%  "
$ws.Range("E10").Value = "  -1.15%  "
$ws.Range("D11").Value = "'41.45"
$ws.Range("E11").Value = "  -0.38%  "
$ws.Range("D12").Value = "'6.192"
$ws.Range("E12").Value = "  -2.28%  "
$ws.Range("D13").Value = "'20.59"
$ws.Range("E13").Value = "  -0.67%  "
$ws.Range("D14").Value = "1.864.98"
$ws.Range("E14").Value = "  -0.62%  "
$ws.Range("D15").Value = "'7.248"
$ws.Range("E15").Value = "  +0.61%  "
$ws.Range("D16").Value = "'1.004"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").Value = "'0.00001097"
$ws.Range("E17").Value = "  -0.80%  "
$ws.Range("D18").Value = "'90.57"
$ws.Range("D20").Value = "'17.68"
$ws.Range("E20").Value = "  -2.48%  "
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("D22").Value = "'6.006"
$ws.Range("E22").Value = "  -1.73%  "
$ws.Range("D23").Value = "28.027.04"
$ws.Range("E23").Value = "  -0.24%  "
$ws.Range("D24").Value = "'11.09"
$ws.Range("E24").Value = "  -2.68%  "
$ws.Range("E25").Value = "  -0.93%  "
$ws.Range("D26").Value = "2.075.06"
$ws.Range("E26").Value = "  -0.87%  "
$ws.Range("D27").Value = "'2.517"
$ws.Range("E27").Value = "  -0.76%  "
$ws.Range("E28").Value = "  +0.74%  "
$ws.Range("D29").Value = "'20.46"
$ws.Range("D30").Value = "'124.59"
$ws.Range("E30").Value = "  -1.68%  "
$ws.Range("D31").Value = "'0.1064"
$ws.Range("E31").Value = "  +1.22%  "
$ws.Range("D32").Value = "'1.028"
$ws.Range("E32").Value = "  -2.91%  "
$ws.Range("D33").Value = "'5.974"
$ws.Range("E33").Value = "  +6.60%  "
$ws.Range("D34").Value = "'3.602"
$ws.Range("E34").Value = "  +0.11%  "
$ws.Range("D35").Value = "'9.352"
$ws.Range("E35").Value = "  -2.71%  "
$ws.Range("D36").Value = "'0.02416"
$ws.Range("E36").Value = "  -0.15%  "
$ws.Range("D37").Value = "'0.06491"
$ws.Range("E37").Value = "  -1.24%  "
$ws.Range("E38").Value = "  -0.38%  "
$ws.Range("D39").Value = "'0.6553"
$ws.Range("E39").Value = "  +2.53%  "
$ws.Range("D40").Value = "'1.195"
$ws.Range("E40").Value = "  -0.94%  "
$ws.Range("D41").Value = "'5.005"
$ws.Range("E41").Value = "  +1.77%  "
$ws.Range("D42").Value = "'1.222"
$ws.Range("E42").Value = "  -4.14%  "
$ws.Range("D43").Value = "'11.15"
$ws.Range("E43").Value = "  -2.82%  "
$ws.Range("D44").Value = "'0.6146"
$ws.Range("E44").Value = "  +1.99%  "
$ws.Range("D45").Value = "'13.06"
$ws.Range("E45").Value = "  -1.06%  "
$ws.Range("D46").Value = "'1.281"
$ws.Range("E46").Value = "  +0.39%  "
$ws.Range("E47").Value = "  -0.46%  "
$ws.Range("D48").Value = "'2.007"
$ws.Range("E49").Value = "  -1.94%  "
$ws.Range("D50").Value = "'120.43"
$ws.Range("E50").Value = "  -0.68%  "
$ws.Range("D51").Value = "'78.40"
$ws.Range("E51").Value = "  -1.57%  "
